$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.551.70"
$ws.Range("E2").Value = "  +2.68%  "
$ws.Range("D3").Value = "3.310.11"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.25"
$ws.Range("E5").Value = "  +5.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.75"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("E8").Value = "  +3.22%  "
$ws.Range("D9").Value = "3.292.80"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.176"
$ws.Range("E10").Value = "  +2.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.576"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.12"
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000273"
$ws.Range("E13").Value = "  +4.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "635.22"
$ws.Range("E14").Value = "  +10.43%  "
$ws.Range("D15").Value = "3.846.25"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.41"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").Value = "67.742.02"
$ws.Range("E17").Value = "  +3.05%  "
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("D19").Value = "3.323.27"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.63"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.87"
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.898"
$ws.Range("E22").Value = "  +1.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.60"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.98"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.99"
$ws.Range("E26").Value = "  +1.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.76"
$ws.Range("E27").Value = "  +3.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.53"
$ws.Range("E28").Value = "  +3.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.61"
$ws.Range("E29").Value = "  +7.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.53"
$ws.Range("E30").Value = "  +1.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.62"
$ws.Range("E31").Value = "  +1.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "590.40"
$ws.Range("E32").Value = "  +6.36%  "
$ws.Range("D33").Value = "3.922.97"
$ws.Range("E33").Value = "  +5.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.90"
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("E35").Value = "  -3.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.103"
$ws.Range("E36").Value = "  +1.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.997"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.54"
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.25"
$ws.Range("E39").Value = "  +4.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.128"
$ws.Range("E40").Value = "  +1.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.67"
$ws.Range("E41").Value = "  +4.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "32.49"
$ws.Range("E42").Value = "  -1.48%  "
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.38"
$ws.Range("E43").Value = "  +1.50%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0679"
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.336"
$ws.Range("E45").Value = "  +2.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0412"
$ws.Range("E46").Value = "  +1.90%  "
$ws.Range("E47").Value = "  +1.89%  "
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.36"
$ws.Range("E50").Value = "  +11.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "130.19"
$ws.Range("E51").Value = "  +4.58%  "

Write-Host "Applied 97 cell updates"
